{"js": "// Remove the final paragraph of the document body \u2014 the gloss\n// \"(19d 29) .i. act ba samlid d\u00faib cid i mm e\u00edcndarcus\" \u2014 leaving the\n// previous paragraph (ending in \"som\") as the new last paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.delete();\nawait context.sync();\n", "ps1": "# Remove the final paragraph of the document body \u2014 the gloss\n# \"(19d 29) .i. act ba samlid d\u00faib cid i mm e\u00edcndarcus\" \u2014 leaving the\n# previous paragraph (ending in \"som\") as the new last paragraph.\n$d = $word.ActiveDocument\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.Delete()\n"}
